# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.416.65'
$ws.Range('E2').Value = '  +2.54%  '
$ws.Range('D3').Value = '2.426.55'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'318.99"
$ws.Range('E5').Value = '  +3.59%  '
$ws.Range('D6').Value = "'103.19"
$ws.Range('E6').Value = '  +3.25%  '
$ws.Range('D7').Value = "'0.517"
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = "'0.530"
$ws.Range('E9').Value = '  +6.08%  '
$ws.Range('D10').Value = "'35.66"
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').Value = "'18.23"
$ws.Range('E13').Value = '  -2.89%  '
$ws.Range('D14').Value = "'7.09"
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').Value = '2.806.03'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').Value = '2.416.42'
$ws.Range('E16').Value = '  -2.68%  '
$ws.Range('D17').Value = "'0.842"
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '45.300.32'
$ws.Range('E18').Value = '  +2.30%  '
$ws.Range('D19').Value = "'12.24"
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').Value = "'6.35"
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('D21').Value = '0.0₃0925'
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').Value = "'68.99"
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').Value = "'245.06"
$ws.Range('E23').Value = '  +1.79%  '
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D27').Value = "'25.81"
$ws.Range('D28').Value = "'2.27"
$ws.Range('E28').Value = '  -3.15%  '
$ws.Range('D29').Value = "'9.62"
$ws.Range('E29').Value = '  +1.04%  '
$ws.Range('D30').Value = "'49.40"
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('D31').Value = "'32.99"
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').Value = "'20.33"
$ws.Range('E32').Value = '  +8.95%  '
$ws.Range('E33').Value = '  +4.06%  '
$ws.Range('D34').Value = "'5.22"
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('D36').Value = "'0.0768"
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = "'1.87"
$ws.Range('E37').Value = '  -2.83%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'4.47"
$ws.Range('E38').Value = '  -2.55%  '
$ws.Range('D39').Value = "'2.88"
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('D40').Value = "'125.66"
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('E42').Value = '  -3.53%  '
$ws.Range('D43').Value = "'20.46"
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').Value = '1.936.74'
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('E46').Value = '  -3.06%  '
$ws.Range('D47').Value = "'2.93"
$ws.Range('E47').Value = '  +1.39%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'1.79"
$ws.Range('E48').Value = '  +8.45%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'9.15"
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('D50').Value = "'76.86"
$ws.Range('E50').Value = '  +4.21%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = "'4.81"
$ws.Range('E51').Value = '  +5.95%  '
